$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EL")
$ws.Name = "RMSE"
